# Apply the commit: add more extras (M8 washers, modified endstop holder
# for limit switches with 24.5mm hole spacing)
#
# Logical change on Sheet1:
#   D5: "1 of 3" -> "q2 (X left only)"
#   D8: "MISSING" -> "q2"
#   Selection moves from D8 to E12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value = "q2"
$ws.Range("D5").Value = "q2 (X left only)"

$ws.Activate()
$ws.Range("E12").Select()
